# Update SAM 2014 beta feedback spreadsheet - Oct 23 entries (rows 57-60)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the cells that introduce NEW shared strings first, in the exact order
# they need to land in the shared-strings table (180 .. 188).
$ws.Range("D57").Value = "Here is the model I would like to simulate.  The model has 208 modules on the flat roof and  70 on the pitched roof.  I modeled this building looks like two but actually it is just one building; so, just one electricity meter.  At least two different string inverters are needed.`nI believe that, for this situation, multiple subsystem modeling feature is much useful.`n"
$ws.Range("D58").Value = ".  I am attaching two screen shots; one is for SAM 2014.1.14  and another is for SAM 2014.9.30.  For PG&E residential, usually it is monthly tier but on new SAM 2014.9.30, there is no monthly tiers."
$ws.Range("H58").Value = "Usability issue with URDB window"
$ws.Range("D59").Value = "The ‘register’ button didn’t look like a button. "
$ws.Range("C59").Value = "Michael F. Troge <mtroge@oneidanation.org>"
$ws.Range("E57").Value = "Followed up. Forwarded to team."
$ws.Range("E58").Value = "Followed up. On meeting agenda."
$ws.Range("D60").Value = "Will  the new version be able to load saved zsam files from the previous version, or TMY.tm2 weather files downloaded in the previous version?"
$ws.Range("C60").Value = "Gomez, Tommaso <tommaso.gomez@intel.com>"

# Remaining cells reusing EXISTING shared strings / plain values.
$ws.Range("B57").Value = "Email"
$ws.Range("C57").Value = "John Yang <john.y@solaronesolution.com>"
$ws.Range("B58").Value = "Email"
$ws.Range("C58").Value = "John Yang <john.y@solaronesolution.com>"
$ws.Range("B59").Value = "Email"
$ws.Range("E59").Value = "Followed up"
$ws.Range("B60").Value = "Email"
$ws.Range("E60").Value = "Followed up"

# Date columns (A = Date, F = Last reply date) -- copy number format from the
# row above so they pick up the existing date style (s="1") instead of General.
$ws.Range("A56").Copy()
$ws.Range("A57").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A57").Value = 41935
$ws.Range("F56").Copy()
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("F57").Value = 41935

$ws.Range("A57").Copy()
$ws.Range("A58").PasteSpecial(-4122)
$ws.Range("A58").Value = 41933
$ws.Range("F57").Copy()
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("F58").Value = 41933

$ws.Range("A58").Copy()
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("A59").Value = 41933
$ws.Range("F58").Copy()
$ws.Range("F59").PasteSpecial(-4122)
$ws.Range("F59").Value = 41933

$ws.Range("A59").Copy()
$ws.Range("A60").PasteSpecial(-4122)
$ws.Range("A60").Value = 41933
$ws.Range("F59").Copy()
$ws.Range("F60").PasteSpecial(-4122)
$ws.Range("F60").Value = 41933

# Row heights to match the wrapped-text content.
$ws.Rows.Item(57).RowHeight = 120
$ws.Rows.Item(58).RowHeight = 45
$ws.Rows.Item(60).RowHeight = 30

# Move the active selection the way the author's last edit left it.
$ws.Range("A61").Select()

Write-Host "Edit complete"
